$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 810 (pushing the existing row 810 "2026/12/29"
# block, and everything after it, down by two rows -> new dimension A1:D853).
$ws.Rows("810:811").Insert()

# New row 810: continuation of 2026/02/14 (Saturday) readings.
$ws.Cells.Item(810, 1).NumberFormat = "@"
$ws.Cells.Item(810, 1).Value = "2026/02/14"
$ws.Cells.Item(810, 2).Value = "土"
$ws.Cells.Item(810, 3).Value = 22
$ws.Cells.Item(810, 4).Value = 201

# New row 811: first reading of the new day, 2026/02/15 (Sunday).
$ws.Cells.Item(811, 1).NumberFormat = "@"
$ws.Cells.Item(811, 1).Value = "2026/02/15"
$ws.Cells.Item(811, 2).Value = "日"
$ws.Cells.Item(811, 3).Value = 0
$ws.Cells.Item(811, 4).Value = 201
